$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like values (e.g. "578.69", "1.00") are written as literal
# text rather than being auto-coerced into numbers (which would also
# strip meaningful trailing zeros). Column D holds price strings as text
# in the source data, so force the Text number format before assigning.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D14", "D15", "D16", "D18", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D29", "D31", "D33", "D36", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.170.11"
$ws.Range("E2").Value = "  +5.11%  "
$ws.Range("D3").Value = "2.759.09"
$ws.Range("E3").Value = "  +4.02%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "578.69"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").Value = "154.52"
$ws.Range("E6").Value = "  +6.34%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").Value = "2.756.90"
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "6.73"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").Value = "3.244.76"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").Value = "26.57"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").Value = "64.052.93"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("E17").Value = "  +5.97%  "
$ws.Range("D18").Value = "2.757.13"
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("D19").Value = "11.97"
$ws.Range("E19").Value = "  +2.19%  "
$ws.Range("D20").Value = "4.88"
$ws.Range("D21").Value = "360.58"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "6.97"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "66.45"
$ws.Range("E25").Value = "  +3.72%  "
$ws.Range("E26").Value = "  +5.14%  "
$ws.Range("D27").Value = "8.55"
$ws.Range("E27").Value = "  +4.50%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "0.0₃0913"
$ws.Range("E29").Value = "  +11.66%  "
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "7.08"
$ws.Range("E31").Value = "  +3.07%  "
$ws.Range("E32").Value = "  +17.69%  "
$ws.Range("D33").Value = "172.16"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("E35").Value = "  +2.44%  "
$ws.Range("D36").Value = "4.82"
$ws.Range("E36").Value = "  +7.17%  "
$ws.Range("E37").Value = "  +8.21%  "
$ws.Range("D38").Value = "1.81"
$ws.Range("E38").Value = "  +9.02%  "
$ws.Range("E39").Value = "  +14.74%  "
$ws.Range("D40").Value = "346.59"
$ws.Range("E40").Value = "  +4.20%  "
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").Value = "39.26"
$ws.Range("E42").Value = "  +1.53%  "
$ws.Range("D43").Value = "5.70"
$ws.Range("E43").Value = "  +9.77%  "
$ws.Range("D44").Value = "21.81"
$ws.Range("E44").Value = "  +6.37%  "
$ws.Range("D45").Value = "21.75"
$ws.Range("E45").Value = "  +6.22%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0589"
$ws.Range("E46").Value = "  +4.29%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.646"
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("D48").Value = "137.04"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("D49").Value = "0.0255"
$ws.Range("E49").Value = "  +2.58%  "
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.25%  "
